$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText([string]$addr, [string]$val) {
    $cell = $ws.Range($addr)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $origStyle
}

Set-CellText "D2" "71.926.99"
Set-CellText "E2" "  -0.37%  "
Set-CellText "D3" "2.683.29"
Set-CellText "E3" "  +1.70%  "
Set-CellText "E4" "  +0.13%  "
Set-CellText "D5" "597.13"
Set-CellText "E5" "  -2.00%  "
Set-CellText "D6" "173.89"
Set-CellText "E6" "  -3.36%  "
Set-CellText "E7" "  +0.09%  "
Set-CellText "E8" "  -0.32%  "
Set-CellText "D9" "2.683.96"
Set-CellText "E9" "  +1.86%  "
Set-CellText "E10" "  -3.18%  "
Set-CellText "E11" "  +1.91%  "
Set-CellText "E12" "  +1.21%  "
Set-CellText "E13" "  -0.85%  "
Set-CellText "D14" "3.174.66"
Set-CellText "E14" "  +2.93%  "
Set-CellText "D15" "71.884.32"
Set-CellText "E15" "  -0.08%  "
Set-CellText "D16" "0.0000183"
Set-CellText "E16" "  -2.71%  "
Set-CellText "D17" "26.14"
Set-CellText "E17" "  -1.89%  "
Set-CellText "D18" "2.684.84"
Set-CellText "E18" "  +2.42%  "
Set-CellText "D19" "12.23"
Set-CellText "E19" "  +5.98%  "
Set-CellText "D20" "8.12"
Set-CellText "E20" "  +0.71%  "
Set-CellText "D21" "370.44"
Set-CellText "E21" "  -3.67%  "
Set-CellText "E22" "  +0.47%  "
Set-CellText "E23" "  -1.07%  "
Set-CellText "D24" "72.19"
Set-CellText "E24" "  -0.99%  "
Set-CellText "E25" "  -0.11%  "
Set-CellText "D26" "4.33"
Set-CellText "E26" "  -3.09%  "
Set-CellText "D27" "9.75"
Set-CellText "E27" "  -2.00%  "
Set-CellText "D28" "2.821.44"
Set-CellText "E28" "  +1.96%  "
Set-CellText "D29" "0.997"
Set-CellText "E29" "  -0.03%  "
Set-CellText "D30" "0.0₃0959"
Set-CellText "E30" "  -0.97%  "
Set-CellText "D31" "8.05"
Set-CellText "E31" "  -0.48%  "
Set-CellText "D32" "497.96"
Set-CellText "E32" "  -9.26%  "
Set-CellText "D33" "1.29"
Set-CellText "E33" "  -3.47%  "
Set-CellText "E34" "  -1.07%  "
Set-CellText "E35" "  -0.01%  "
Set-CellText "D36" "164.17"
Set-CellText "E36" "  -1.47%  "
Set-CellText "D37" "19.57"
Set-CellText "E37" "  +1.40%  "
Set-CellText "D38" "19.09"
Set-CellText "E38" "  -0.26%  "
Set-CellText "E39" "  -1.89%  "
Set-CellText "E40" "  -6.57%  "
Set-CellText "E41" "  -4.79%  "
Set-CellText "E42" "  +0.02%  "
Set-CellText "D43" "5.00"
Set-CellText "E43" "  -0.97%  "
Set-CellText "E44" "  -0.19%  "
Set-CellText "D45" "2.54"
Set-CellText "E45" "  -2.98%  "
Set-CellText "D46" "157.38"
Set-CellText "E46" "  +4.27%  "
Set-CellText "D47" "39.34"
Set-CellText "E47" "  -0.61%  "
Set-CellText "E48" "  +4.77%  "
Set-CellText "D49" "3.73"
Set-CellText "E49" "  +1.51%  "
Set-CellText "E50" "  +3.42%  "
Set-CellText "D51" "0.0761"
Set-CellText "E51" "  +0.59%  "
